# Auto-generated script to append new log rows to the "Logs" sheet
# matching the target diff (Thursday Nov 3 / Friday Nov 4 / Monday Nov 7, 2016 entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Seed the 6 brand-new shared strings in the exact order they are first used,
# using a disposable scratch sheet so the shared-string table gets indices 274-279
# assigned in the same order as in the target workbook. ---
$seed = $wb.Worksheets.Add()
$seed.Cells.Item(1,1).Value = "Built in PC not working - demo roll in PC cart in room already - make sure client okay."
$seed.Cells.Item(2,1).Value = "Pick up roll in PC cart from CLH C - return sound cable and VGA cable to drawer of podium. Return roll in PC cart to Lassonde 1011 storeroom. "
$seed.Cells.Item(3,1).Value = "LEAVE PORTABLE SCREEN IN ROOM"
$seed.Cells.Item(4,1).Value = "Turn on SMALL PA SYSTEM - don't use 5065*0 - won't work today. Use Small Pa speaker for audio."
$seed.Cells.Item(5,1).Value = "Pick up SMALL PA SPEAKER ON CART. Return to Lassonde 1011 storeoroom."
$seed.Cells.Item(6,1).Value = "Press `"MUTE`" button on wireless mics to turn on mics."
$excel.DisplayAlerts = $false
$seed.Delete()
$excel.DisplayAlerts = $true

# --- Template rows already present in the sheet, used as format sources ---
$sectionTemplate = $ws.Range("A5:F5")   # blank day-header row style (s=21,22,23,21,24,25)
$dataTemplate = $ws.Range("A6:F6")       # standard data row style (s=3,10,12,9,11,19)

function Set-SectionRow($rowNum, $dayName) {
    $dst = $ws.Range("A" + $rowNum + ":F" + $rowNum)
    $sectionTemplate.Copy()
    $dst.PasteSpecial(-4122)
    $ws.Cells.Item($rowNum, 2).Value = $dayName
}

function Set-DataRow($rowNum, $a, $b, $c, $d, $e, $f, $ht, $boldF) {
    $dst = $ws.Range("A" + $rowNum + ":F" + $rowNum)
    $dataTemplate.Copy()
    $dst.PasteSpecial(-4122)
    $ws.Cells.Item($rowNum, 1).Value = $a
    $ws.Cells.Item($rowNum, 2).Value = $b
    $ws.Cells.Item($rowNum, 3).Value = $c
    $ws.Cells.Item($rowNum, 4).Value = $d
    $ws.Cells.Item($rowNum, 5).Value = $e
    $ws.Cells.Item($rowNum, 6).Value = $f
    if ($boldF) {
        $ws.Cells.Item($rowNum, 6).Font.Bold = $true
    }
    if ($ht) {
        $ws.Rows.Item($rowNum).RowHeight = $ht
    }
}


Set-SectionRow 637 "THURSDAY"
Set-DataRow 638 "AV Shutdown" 42677 "1900" "CLH" "L" "PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN'T GET LOGGED OFF." 45 $False
Set-DataRow 639 "AV Shutdown" 42677 "1730" "LSB" "101" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 640 "AV Shutdown" 42677 "1900" "LSB" "103" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 641 "AV Shutdown" 42677 "1900" "LSB" "105" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 642 "AV Shutdown" 42677 "1730" "LSB" "107" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 643 "Demo" 42677 "1730" "CLH" "C" "Built in PC not working - demo roll in PC cart in room already - make sure client okay." 30 $False
Set-DataRow 644 "Pickup PC" 42677 "1900" "CLH" "C" "Pick up roll in PC cart from CLH C - return sound cable and VGA cable to drawer of podium. Return roll in PC cart to Lassonde 1011 storeroom. " 45 $False
Set-DataRow 645 "Lockup" 42677 "2030" "CLH" "K" "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS." 30 $False
Set-DataRow 646 "Lockup" 42677 "2150" "CLH" "M" "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS." 30 $False
Set-DataRow 647 "Lockup" 42677 "2150" "CLH" "H" "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS." 30 $False
Set-DataRow 648 "Lockup" 42677 "2150" "CLH" "J" "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS." 30 $False
Set-SectionRow 653 "FRIDAY"
Set-DataRow 654 "Pickup PC" 42678 "1630" "SC" "MDR" "Pick up roll in PC and Projector carts, all matts and cables and return equipment to Bethune 201 storeroom. Key for Stong MDR is in CB 121A storeroom." 45 $False
Set-DataRow 655 "Other" 42678 "1630" "SC" "MDR" "LEAVE PORTABLE SCREEN IN ROOM" 0 $False
Set-DataRow 656 "AV Shutdown" 42678 "1730" "CLH" "L" "PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN'T GET LOGGED OFF." 45 $False
Set-DataRow 657 "Lockup" 42678 "1730" "CLH" "H" "PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS." 30 $False
Set-SectionRow 662 "MONDAY"
Set-DataRow 663 "AV Shutdown" 42681 "1530" "BC" "320" "Please turn off PC and Projector. Leave projector remote on PC cart and lock room.  EQUIPMENT STAYS IN ROOM. Key for room in CB 121A on Bethune classroom keys." 45 $False
Set-DataRow 664 "Setup Mic" 42681 "1715" "LAS" "C" "Take cart with mixer, 2 wireless mics and 2 mic stands from Lassonde 1011 storeroom (across from Lassonde A). Go to Lassonde C classroom (class starts at 5:30 pm but be there early in case previous class ends early). " 60 $False
Set-DataRow 665 "Other" 42681 "1715" "LAS" "C" "Turn on SMALL PA SYSTEM - don't use 5065*0 - won't work today. Use Small Pa speaker for audio." 30 $True
Set-DataRow 666 "Other" 42681 "1715" "LAS" "C" "Plug in power cord from cart on to power outlet on left side of podium (to left of document camera). Turn on mixer. Turn on wireless microphone receivers on cart (NOTE: DO NOT PRESS `"SYNC`" BUTTON`" - POWER BUTTON IS FIRST BUTTON TO THE RIGHT ON RECEIVER). " 75 $False
Set-DataRow 667 "Other" 42681 "1715" "LAS" "C" "Press `"MUTE`" button on wireless mics to turn on mics." 0 $False
Set-DataRow 668 "Other" 42681 "1715" "LAS" "C" "Once volumes are set, place one mic stand with mic halfway up aisle on right and one mic stand with mic halfway up aisle on left. Demo volume controls to prof. and demo PC. Leave microphone bags with milk carton on cart in room. PLEASE FIND OUT END TIME OF CLASS FROM PROF. AND TELL MASI AS MICROPHONES ARE EXPENSIVE. TELL PROF. TO STAY WITH MICS UNTIL THEY ARE PICKED UP. TELL HIM TO CALL ext 55800   WHEN DONE (use phone in classroom)." 120 $False
Set-DataRow 669 "AV Shutdown" 42681 "1900" "CLH" "L" "PLEASE MAKE SURE CRESTRON GETS LOGGED OFF. WE ARE HAVING PROBLEMS WITH THIS ROOM WHEN IT DOESN'T GET LOGGED OFF." 45 $False
Set-DataRow 670 "AV Shutdown" 42681 "1630" "LSB" "101" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 671 "AV Shutdown" 42681 "1900" "LSB" "103" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 672 "AV Shutdown" 42681 "1900" "LSB" "105" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 673 "AV Shutdown" 42681 "1730" "LSB" "106" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 674 "AV Shutdown" 42681 "1900" "LSB" "107" "Make sure neck mic goes back to drawer and log off touchscreen." 0 $False
Set-DataRow 675 "Pickup Mic" 42681 "1850" "LAS" "C" "Pick up 2 wireless mics on stands with cart. Move all equipment on cart - cart has 2 wireless mic receivers and mixer and mic cables. Pick up 2 mic stands - return all equipment to Lassonde 1011 storeroom (across the hall from Lassonde A). PLEASE PUT 2 WIRELESS MICS IN BAGS PROVIDED IN MILK CARTON ON CART. Very expensive mics - please go early and treat mics with care." 90 $False
Set-DataRow 676 "Other" 42681 "1850" "LAS" "C" "Pick up SMALL PA SPEAKER ON CART. Return to Lassonde 1011 storeoroom." 30 $True
Set-DataRow 677 "Other" 42681 "1850" "LAS" "C" "Turn off wireless microphones by pressing `"MUTE`" button on mics." 30 $False
Set-DataRow 678 "Other" 42681 "1850" "LAS" "C" "Turn off wireless microphone receivers by pressing `"POWER`" button and not `"SYNC`" button. " 30 $False
Set-DataRow 679 "Other" 42681 "1850" "LAS" "C" "PLEASE BE ON TIME - Prof upset last week when no one came till 7:05 pm and other class was starting." 30 $False

# --- Update the active selection to match the new end of the log ---
$ws.Range("F682").Select()
try { $excel.ActiveWindow.ScrollRow = 668 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}

Write-Output "rows added"